$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("NSAA", "position", "dhc", 60),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "AD", "acts", 10),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "AD", "acts", 10)
)

$startRow = 277
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
